$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range to find the last row with data in column J
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Apply a custom date/time number format to column J cells (rows 2..lastRow)
$dateFormat = "YYYY-MM-DD HH:MM:SS"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
    $cell.NumberFormat = $dateFormat
    $cell.Value = "2023-03-03"
}
